$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - __init__ / Attribute set to input values
$ws.Range("F7").Value = '(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("G7").Value = "Pass/No errors"

# Row 8 - __init__ / Exception raised when title is blank
$ws.Range("F8").Value = '(1, "", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("G8").Value = "Raise ValueError"

# Row 9 - __init__ / Exception raised when author is blank
$ws.Range("F9").Value = '(1, "DUNE", "", Genre.FICTION, False)'
$ws.Range("G9").Value = "Raise ValueError"

# Row 10 - __init__ / Exception raised when invalid Genre
$ws.Range("F10").Value = '(1, "DUNE", "Frank Herbert", "INVALID, False)'
$ws.Range("G10").Value = "Raise ValueError"

# Row 11 - title / returns title attribute
$ws.Range("E11").Value = 'LibraryItem(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("G11").Value = "Returns title"

# Row 12 - author / returns author attribute
$ws.Range("E12").Value = 'LibraryItem(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("G12").Value = "Returns author"

# Row 13 - Genre / returns Genre attribute
$ws.Range("E13").Value = 'LibraryItem(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("G13").Value = "Returns genre"

# Row 14 - new: __init__ / Exception raised when invalid item id
$ws.Range("C14").Value = "__init__"
$ws.Range("D14").Value = "Exception raised when invalid item id"
$ws.Range("E14").Value = "None"
$ws.Range("F14").Value = '("INVALID ID", "DUNE",  "Frank Herbert", "INVALID", False)'
$ws.Range("G14").Value = "Raises ValueError"

# Row 15 - new: __init__ / exception raised when invalid is borrowed
$ws.Range("C15").Value = "__init__"
$ws.Range("D15").Value = "exception raised when invalid is borrowed"
$ws.Range("E15").Value = "None"
$ws.Range("F15").Value = '(1, "DUNE",  "Frank Herbert", "INVALID", "Not a boolean")'
$ws.Range("G15").Value = "Raises ValueError"

# Row 16 - new: item_id / returns item_id attribute
$ws.Range("C16").Value = "item_id"
$ws.Range("D16").Value = "returns item_id attribute"
$ws.Range("E16").Value = 'LibraryItem(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "Returns item_Id"

# Row 17 - new: is_ borrowed / return is_borrowed attribute
$ws.Range("C17").Value = "is_ borrowed"
$ws.Range("D17").Value = "return is_borrowed attribute"
$ws.Range("E17").Value = 'LibraryItem(1, "DUNE", "Frank Herbert", Genre.FICTION, False)'
$ws.Range("F17").Value = "None"
$ws.Range("G17").Value = "returns is_borrowed"

# Update sheet view: scroll position and active selection
$excel.Goto($ws.Range("A8"), $true)
$ws.Range("J17").Select()
